$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A (codes) and C (prices, e.g. "$ 62.990") look numeric to Excel's
# type inference and would otherwise be auto-converted to numbers/currency.
# Force them to Text format first so the literal strings from the diff are
# preserved. Column B ("Disponible") is non-numeric text already, and column
# D holds genuine numeric quantities, so neither needs a format change.
$ws.Range("A2:A5").NumberFormat = "@"
$ws.Range("C2:C5").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "10962736022"
$ws.Range("B2").Value = "Disponible"
$ws.Range("C2").Value = "$ 62.990"
$ws.Range("D2").Value = 3

# Row 3
$ws.Range("A3").Value = "10962389016"
$ws.Range("B3").Value = "Disponible"
$ws.Range("C3").Value = "$ 165.990"
$ws.Range("D3").Value = 3

# Row 4
$ws.Range("A4").Value = "10962389018"
$ws.Range("B4").Value = "Disponible"
$ws.Range("C4").Value = "$ 165.990"
$ws.Range("D4").Value = 3

# Row 5
$ws.Range("A5").Value = "10930745010"
$ws.Range("B5").Value = "Disponible"
$ws.Range("C5").Value = "$ 182.990"
$ws.Range("D5").Value = 11
